$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.837.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5017"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06400"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07688"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.650.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.859.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5425"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7917"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.863.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.318"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.926"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.966"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.915"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1139"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.695"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04981"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.254"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.85%  "
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.536"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.166.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8918"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.617"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5601"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01555"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.549"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8067"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.771.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈115"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.002"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05080"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "

Write-Host "Applied cryptos update"
